$d = $word.ActiveDocument

# Anchor on the "SMARTREWARDS FAQ'S" heading paragraph - keep it, but
# remove every FAQ question/answer paragraph that follows it, down
# through (and including) the final answer paragraph ending in
# "...will get paid." The trailing empty paragraph after that (the very
# last paragraph of the document) is left untouched.

$headingRange = $d.Content
[void]$headingRange.Find.Execute("SMARTREWARDS FAQ")
[void]$headingRange.Expand(4)   # wdParagraph - expand to the whole heading paragraph (incl. mark)

$lastAnswerRange = $d.Content
[void]$lastAnswerRange.Find.Execute("will get paid.")
[void]$lastAnswerRange.Expand(4)   # wdParagraph - expand to the whole last FAQ answer paragraph

$toDelete = $d.Range($headingRange.End, $lastAnswerRange.End)
$toDelete.Delete()
